$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Step 1: row 3 (08:45 - 09:00 | COLFRESH COFFEE) -> time becomes 09:00 - 09:15
$t.Cell(3, 1).Range.Text = "09:00 - 09:15"

# --- Step 2: row 4 (09:00 - 09:15 | ARMANDO VELASQUEZ) -> time becomes 09:15 - 09:30,
#     buyer becomes FLOR A FRUTO (the slot that used to belong to the now-removed
#     10:00 - 10:15 row)
$t.Cell(4, 1).Range.Text = "09:15 - 09:30"
$t.Cell(4, 3).Range.Text = "FLOR A FRUTO"

# --- Step 3: insert a brand-new row right before the old "10:00 - 10:15 | FLOR A FRUTO"
#     row (currently row 5), carrying the 09:30 - 09:45 | ARMANDO VELASQUEZ appointment
#     that moved out of row 4.
$refRow = $t.Rows.Item(5)
$newRow = $t.Rows.Add($refRow)
$newRow.Cells.Item(1).Range.Text = "09:30 - 09:45"
$newRow.Cells.Item(3).Range.Text = "ARMANDO VELÁSQUEZ"

# --- Step 4: the old "10:00 - 10:15 | FLOR A FRUTO" row (now pushed down to row 6)
#     is obsolete -- its content already lives in the new rows above -- so delete it.
$t.Rows.Item(6).Delete()

# --- Step 5: last row (11:15 - 11:30 | NEIRA YORK COFFEE) -> time becomes 11:30 - 11:45
$t.Cell($t.Rows.Count, 1).Range.Text = "11:30 - 11:45"
